# Update the date/title line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-10-11 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-12 Sunday", 2)

# Update each division-fact cell in the single table by explicit (row, column)
# coordinates, so that duplicate values (e.g. multiple "53÷3=17, 2" cells)
# are each replaced with their own correct new value instead of a blanket
# find/replace that would overwrite every matching occurrence identically.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "18÷6=3, 0"
$tbl.Cell(1, 2).Range.Text = "25÷7=3, 4"
$tbl.Cell(1, 3).Range.Text = "17÷2=8, 1"
$tbl.Cell(1, 4).Range.Text = "30÷3=10, 0"
$tbl.Cell(1, 5).Range.Text = "37÷7=5, 2"

$tbl.Cell(5, 1).Range.Text = "58÷2=29, 0"
$tbl.Cell(5, 2).Range.Text = "17÷8=2, 1"
$tbl.Cell(5, 3).Range.Text = "13÷4=3, 1"
$tbl.Cell(5, 4).Range.Text = "42÷4=10, 2"
$tbl.Cell(5, 5).Range.Text = "91÷2=45, 1"

$tbl.Cell(9, 1).Range.Text = "72÷3=24, 0"
$tbl.Cell(9, 2).Range.Text = "23÷9=2, 5"
$tbl.Cell(9, 3).Range.Text = "90÷2=45, 0"
$tbl.Cell(9, 4).Range.Text = "24÷6=4, 0"
$tbl.Cell(9, 5).Range.Text = "98÷8=12, 2"

$tbl.Cell(13, 1).Range.Text = "68÷4=17, 0"
$tbl.Cell(13, 2).Range.Text = "47÷4=11, 3"
$tbl.Cell(13, 3).Range.Text = "46÷2=23, 0"
$tbl.Cell(13, 4).Range.Text = "35÷4=8, 3"
$tbl.Cell(13, 5).Range.Text = "61÷5=12, 1"

$tbl.Cell(17, 1).Range.Text = "71÷9=7, 8"
$tbl.Cell(17, 2).Range.Text = "21÷4=5, 1"
$tbl.Cell(17, 3).Range.Text = "59÷3=19, 2"
$tbl.Cell(17, 4).Range.Text = "12÷4=3, 0"
$tbl.Cell(17, 5).Range.Text = "95÷5=19, 0"
